# Update cryptos worksheet: refresh price/volume data (and a few
# coin rows that got reordered) per the Fri Mar 10 09:44:48 UTC 2023 feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as TEXT (matches the source data,
# which stores prices as plain strings) without leaving the cell's style
# changed -- we flip to a text format just long enough to assign the
# value, then clear the format back off again.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "19.906.37"
$ws.Range("E2").Value = "  -8.13%  "

$ws.Range("D3").Value = "1.403.41"
$ws.Range("E3").Value = "  -8.46%  "

Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue $ws.Range("D5") "1.002"
$ws.Range("E5").Value = "  +0.06%  "

Set-TextValue $ws.Range("D6") "269.56"
$ws.Range("E6").Value = "  -6.51%  "

Set-TextValue $ws.Range("D7") "0.3668"
$ws.Range("E7").Value = "  -6.67%  "

Set-TextValue $ws.Range("D8") "0.3037"
$ws.Range("E8").Value = "  -3.83%  "

Set-TextValue $ws.Range("D9") "39.13"
$ws.Range("E9").Value = "  -7.50%  "

Set-TextValue $ws.Range("D10") "0.06451"
$ws.Range("E10").Value = "  -9.88%  "

Set-TextValue $ws.Range("D11") "0.9701"
$ws.Range("E11").Value = "  -7.10%  "

Set-TextValue $ws.Range("D12") "1.003"
$ws.Range("E12").Value = "  +0.21%  "

Set-TextValue $ws.Range("D13") "5.266"
$ws.Range("E13").Value = "  -6.30%  "

Set-TextValue $ws.Range("D14") "6.051"
$ws.Range("E14").Value = "  -8.21%  "

Set-TextValue $ws.Range("D15") "16.57"
$ws.Range("E15").Value = "  -10.17%  "

$ws.Range("D16").Value = "1.406.44"
$ws.Range("E16").Value = "  -8.45%  "

Set-TextValue $ws.Range("D17") "0.000009998"
$ws.Range("E17").Value = "  -8.38%  "

Set-TextValue $ws.Range("D18") "0.05671"
$ws.Range("E18").Value = "  -13.96%  "

Set-TextValue $ws.Range("D19") "71.60"
$ws.Range("E19").Value = "  -13.65%  "

Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  +0.08%  "

Set-TextValue $ws.Range("D21") "5.483"
$ws.Range("E21").Value = "  -9.94%  "

Set-TextValue $ws.Range("D22") "14.12"
$ws.Range("E22").Value = "  -8.21%  "

Set-TextValue $ws.Range("D23") "10.56"
$ws.Range("E23").Value = "  -2.21%  "

Set-TextValue $ws.Range("D24") "2.269"
$ws.Range("E24").Value = "  -4.97%  "

$ws.Range("D25").Value = "19.919.18"
$ws.Range("E25").Value = "  -8.10%  "

Set-TextValue $ws.Range("D26") "2.199"
$ws.Range("E26").Value = "  -5.47%  "

Set-TextValue $ws.Range("D27") "135.82"
$ws.Range("E27").Value = "  -7.94%  "

Set-TextValue $ws.Range("D28") "16.51"
$ws.Range("E28").Value = "  -9.68%  "

$ws.Range("D29").Value = "1.566.47"
$ws.Range("E29").Value = "  -8.58%  "

Set-TextValue $ws.Range("D30") "106.96"
$ws.Range("E30").Value = "  -8.45%  "

Set-TextValue $ws.Range("D31") "3.857"
$ws.Range("E31").Value = "  -20.38%  "

Set-TextValue $ws.Range("D32") "5.183"
$ws.Range("E32").Value = "  -11.23%  "

Set-TextValue $ws.Range("D33") "0.7961"
$ws.Range("E33").Value = "  -16.39%  "

Set-TextValue $ws.Range("D34") "0.07617"
$ws.Range("E34").Value = "  -6.11%  "

Set-TextValue $ws.Range("D35") "8.300"
$ws.Range("E35").Value = "  -2.65%  "

Set-TextValue $ws.Range("D36") "0.05737"
$ws.Range("E36").Value = "  -5.23%  "

Set-TextValue $ws.Range("D38") "4.683"
$ws.Range("E38").Value = "  -7.79%  "

Set-TextValue $ws.Range("D39") "0.1906"
$ws.Range("E39").Value = "  -5.46%  "

Set-TextValue $ws.Range("D40") "0.02008"
$ws.Range("E40").Value = "  -8.77%  "

$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D41") "1.340"
$ws.Range("E41").Value = "  -7.42%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "10.08"
$ws.Range("E42").Value = "  -6.71%  "

Set-TextValue $ws.Range("D43") "1.052"
$ws.Range("E43").Value = "  -10.62%  "

Set-TextValue $ws.Range("D44") "0.5220"
$ws.Range("E44").Value = "  -8.81%  "

Set-TextValue $ws.Range("D45") "3.486"
$ws.Range("E45").Value = "  -6.25%  "

Set-TextValue $ws.Range("D46") "11.90"
$ws.Range("E46").Value = "  -8.40%  "

Set-TextValue $ws.Range("D47") "0.5018"
$ws.Range("E47").Value = "  -8.33%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "108.99"
$ws.Range("E48").Value = "  -5.39%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.752"
$ws.Range("E49").Value = "  -5.76%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D50") "1.001"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D51") "1.026"
$ws.Range("E51").Value = "  -11.03%  "
